# Adds columns I ("I0") and J ("IF") with per-row numeric data.
# Header cells I1/J1 copy the format (bold + border + centered) from H1,
# matching the existing header style used by B1:H1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: copy H1's format onto I1:J1, then set the header text ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows: (row, I-value, J-value) ---
$ijData = @(
  @(2, 3, 5),
  @(3, 6, 7),
  @(4, 6, 6),
  @(5, 8, 8),
  @(6, 7, 7),
  @(7, 8, 8),
  @(8, 9, 9),
  @(9, 11, 11),
  @(10, 8, 8),
  @(11, 8, 8),
  @(12, 7, 7),
  @(13, 8, 8),
  @(14, 8, 8),
  @(15, 8, 8),
  @(16, 8, 8),
  @(17, 8, 8),
  @(18, 7, 7),
  @(19, 7, 8),
  @(20, 11, 11),
  @(21, 6, 6),
  @(22, 7, 7),
  @(23, 9, 10),
  @(24, 8, 8),
  @(25, 8, 8),
  @(26, 8, 8),
  @(27, 8, 8),
  @(28, 8, 9),
  @(29, 7, 7),
  @(30, 7, 7),
  @(31, 8, 8),
  @(32, 10, 10),
  @(33, 8, 8),
  @(34, 8, 8),
  @(35, 6, 6),
  @(36, 7, 8),
  @(37, 9, 9),
  @(38, 7, 8),
  @(39, 7, 8),
  @(40, 7, 7),
  @(41, 5, 6),
  @(42, 7, 7),
  @(43, 8, 8),
  @(44, 8, 8),
  @(45, 7, 7),
  @(46, 10, 10),
  @(47, 8, 9),
  @(48, 7, 7),
  @(49, 8, 8),
  @(50, 8, 8),
  @(51, 9, 9),
  @(52, 6, 6),
  @(53, 8, 8),
  @(54, 9, 9),
  @(55, 7, 7),
  @(56, 8, 8),
  @(57, 8, 8),
  @(58, 7, 8),
  @(59, 11, 11),
  @(60, 8, 8),
  @(61, 8, 8),
  @(62, 9, 9),
  @(63, 7, 8),
  @(64, 7, 7),
  @(65, 7, 8),
  @(66, 8, 8),
  @(67, 6, 7),
  @(68, 6, 6),
  @(69, 8, 8),
  @(70, 4, 4),
  @(71, 8, 8),
  @(72, 5, 6),
  @(73, 8, 8),
  @(74, 4, 4),
  @(75, 3, 3)
)

foreach ($entry in $ijData) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}
